$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): version labels - mark older versions as archived,
#    and flag v1.4 as the current version.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "v0 (archived)"
$ws.Range("C1").Value = "v1 (archived)"
$ws.Range("D1").Value = "v1.1 (archived)"
$ws.Range("E1").Value = "v1.2 (archived)"
$ws.Range("F1").Value = "v1.3 (archived)"
$ws.Range("G1").Value = "v1.4 - CURRENT"

# ---------------------------------------------------------------------------
# 2) Row 5 ("Update - functional/structural"): refresh the v1.3 and v1.4
#    short-term-goal notes.
# ---------------------------------------------------------------------------
$ws.Range("G5").Value = "Streamlined dedup blocking code`nFixed bugs from upgrade to pandas 2`nFixed bugs in 1:M match rate calculation"
$ws.Range("H5").Value = "Update workflow to have a csv path and a postgres path`nClean up pipeline to have better integration with the preprocessing component`nMore user-friendly config (yaml files, split preprocessing, etc.)`nHandle matches with multiple IDs, potentially with different strengths`nDraft and pilot match analysis to include in linkage memo;"

# ---------------------------------------------------------------------------
# 3) Row 6 ("Update - documentation/other"): refresh the v1.4 note.
# ---------------------------------------------------------------------------
$ws.Range("H6").Value = "Draft user manual;`nPlanning for case studies of sensitivity testing"

# ---------------------------------------------------------------------------
# 4) "Link to logic" row (row 9): the v1 / v1.1 design-doc hyperlinks are
#    retired; v1.2-v1.4 all point to "Same as v1.1".
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("G9").Value = "Same as v1.1"

# ---------------------------------------------------------------------------
# 5) "Link to commit" row (row 10): retire the v0 commit-link hyperlinks.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""

# ---------------------------------------------------------------------------
# 6) View state: scroll the frozen pane so D5 is the top-left visible cell
#    and F5 is the active selection.
# ---------------------------------------------------------------------------
$ws.Range("F5").Select()
